$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data updates: toggle the Yes/No "Started" flag (column C) for the
#     rows whose status changed in this revision ---
$ws.Range("C8").Value  = "No"
$ws.Range("C9").Value  = "Yes"
$ws.Range("C18").Value = "No"
$ws.Range("C20").Value = "Yes"
$ws.Range("C21").Value = "Yes"
$ws.Range("C22").Value = "No"
$ws.Range("C32").Value = "No"
$ws.Range("C33").Value = "Yes"
$ws.Range("C56").Value = "Yes"
$ws.Range("C60").Value = "No"
$ws.Range("C76").Value = "No"
$ws.Range("C79").Value = "Yes"
$ws.Range("C82").Value = "No"
$ws.Range("C83").Value = "Yes"

# --- View state: move the frozen-pane scroll position / selection to
#     match where the author left off editing ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 63
$ws.Range("C84").Select()
